$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "323.06";  "E2"  = "-2.59%"
    "D3"  = "42.93";   "E3"  = "-5.53%"
    "D4"  = "5.183";   "E4"  = "-7.68%"
    "D5"  = "0.08196"; "E5"  = "-1.91%"
    "E6"  = "-2.69%"
    "D7"  = "1.833";   "E7"  = "-11.61%"
    "D8"  = "0.9332";  "E8"  = "-3.01%"
    "D9"  = "0.1114";  "E9"  = "-4.88%"
    "D10" = "0.1868";  "E10" = "-2.66%"
    "D11" = "0.09493"; "E11" = "-3.67%"
    "D12" = "0.04625"; "E12" = "0.10%"
    "D13" = "7.423";   "E13" = "-28.66%"
    "D14" = "0.1057";  "E14" = "-0.31%"
    "D15" = "0.001294";"E15" = "1.49%"
    "D16" = "0.005852";"E16" = "-4.33%"
    "D17" = "3.364";   "E17" = "-0.40%"
    "D19" = "0.3374";  "E19" = "0.98%"
    "D20" = "0.1388";  "E20" = "-0.43%"
    "D21" = "0.2492";  "E21" = "-6.18%"
    "D22" = "0.04159"; "E22" = "-0.71%"
    "D23" = "0.001244";"E23" = "-5.40%"
    "D24" = "0.004345";"E24" = "-4.90%"
    "E25" = "-8.00%"
    "D26" = "0.0002977";"E26" = "-20.62%"
    "D38" = "0.02757"; "E38" = "1.89%"
    "D39" = "0.05593"; "E39" = "-2.89%"
    "D40" = "0.008313";"E40" = "6.25%"
    "D41" = "0.1396";  "E41" = "-2.60%"
    "D42" = "0.006535";"E42" = "-10.21%"
    "E43" = "3.73%"
    "D44" = "0.007530";"E44" = "-17.61%"
    "D45" = "0.3502";  "E45" = "-1.16%"
    "D46" = "0.00006995";"E46" = "-1.80%"
    "E47" = "-0.34%"
    "D48" = "0.003478";"E48" = "-0.42%"
    "D49" = "0.003528";"E49" = "0.58%"
    "D50" = "0.00002099";"E50" = "-0.34%"
    "D51" = "0.0001999";"E51" = "-0.34%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
